$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2630102.5
$ws.Range("J17").Value = 2673934.8
$ws.Range("L17").Value = 8021804.399999999
$ws.Range("N17").Value = -8022140.399999999

$ws.Range("H53").Value = 63416.125
$ws.Range("I53").Value = 144422.86
$ws.Range("J53").Value = 410.8889
$ws.Range("K53").Value = 144422.86
$ws.Range("L53").Value = 410.8889
$ws.Range("M53").Value = -143785.86
$ws.Range("N53").Value = -1684.8889

$ws.Range("H132").Value = 655.4583
$ws.Range("I132").Value = 605.57446
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 1816.72338
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = 713.2766199999999
$ws.Range("N132").Value = -14060

$ws.Range("H137").Value = 973.39685
$ws.Range("I137").Value = 849.4386
$ws.Range("J137").Value = 2151
$ws.Range("K137").Value = 2548.3158
$ws.Range("L137").Value = 6453
$ws.Range("M137").Value = 1.684200000000146
$ws.Range("N137").Value = -11553

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5172.1772
$ws.Range("I32").Value = 3931.3655
$ws.Range("J32").Value = 11624.4
$ws.Range("K32").Value = 3931.3655
$ws.Range("L32").Value = 11624.4
$ws.Range("M32").Value = -3644.3655
$ws.Range("N32").Value = -12198.4

$ws.Range("H61").Value = 5208.2144
$ws.Range("I61").Value = 5675.4165
$ws.Range("J61").Value = 2405
$ws.Range("K61").Value = 5675.4165
$ws.Range("L61").Value = 2405
$ws.Range("M61").Value = -5463.4165
$ws.Range("N61").Value = -2829

$ws.Range("H74").Value = 1606.5652
$ws.Range("I74").Value = 1413.9445
$ws.Range("K74").Value = 1413.9445
$ws.Range("M74").Value = -539.9445000000001

$ws.Range("H77").Value = 1606.5652
$ws.Range("I77").Value = 1413.9445
$ws.Range("K77").Value = 7069.7225
$ws.Range("M77").Value = -2701.7225

$ws.Range("H97").Value = 776
$ws.Range("I97").Value = 501.1111
$ws.Range("J97").Value = 1023.4
$ws.Range("K97").Value = 501.1111
$ws.Range("L97").Value = 1023.4
$ws.Range("M97").Value = -5.111100000000022
$ws.Range("N97").Value = -2015.4

$ws.Range("H122").Value = 952024.0600000001
$ws.Range("I122").Value = 1168037.4
$ws.Range("J122").Value = 1565.2
$ws.Range("K122").Value = 3504112.2
$ws.Range("L122").Value = 4695.6
$ws.Range("M122").Value = -3501662.2
$ws.Range("N122").Value = -9595.6

$ws.Range("H132").Value = 2587.4897
$ws.Range("I132").Value = 1426.6471
$ws.Range("J132").Value = 5218.7334
$ws.Range("K132").Value = 4279.9413
$ws.Range("L132").Value = 15656.2002
$ws.Range("M132").Value = -1749.9413
$ws.Range("N132").Value = -20716.2002

$ws.Range("H136").Value = 5208.2144
$ws.Range("I136").Value = 5675.4165
$ws.Range("J136").Value = 2405
$ws.Range("K136").Value = 17026.2495
$ws.Range("L136").Value = 7215
$ws.Range("M136").Value = -14476.2495
$ws.Range("N136").Value = -12315

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 42780
$ws.Range("J132").Value = 42780
$ws.Range("L132").Value = 42780
$ws.Range("N132").Value = -52900

$ws.Range("H134").Value = 4033.1667
$ws.Range("I134").Value = 4305.6763
$ws.Range("K134").Value = 12917.0289
$ws.Range("M134").Value = -10382.0289

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3158.9302
$ws.Range("I31").Value = 1435.5454
$ws.Range("K31").Value = 1435.5454
$ws.Range("M31").Value = -1140.5454

$ws.Range("H34").Value = 3158.9302
$ws.Range("I34").Value = 1435.5454
$ws.Range("K34").Value = 1435.5454
$ws.Range("M34").Value = -1233.5454

$ws.Range("H58").Value = 1371.1578
$ws.Range("I58").Value = 1012.0833
$ws.Range("K58").Value = 1012.0833
$ws.Range("M58").Value = -809.0833

$ws.Range("H107").Value = 290.6207
$ws.Range("I107").Value = 262.16666
$ws.Range("K107").Value = 262.16666
$ws.Range("M107").Value = 1657.83334

$ws.Range("H132").Value = 1954.4286
$ws.Range("I132").Value = 1776.9429
$ws.Range("J132").Value = 2398.1428
$ws.Range("K132").Value = 5330.8287
$ws.Range("L132").Value = 7194.428400000001
$ws.Range("M132").Value = -2800.8287
$ws.Range("N132").Value = -12254.4284

$ws.Range("H134").Value = 2664.0344
$ws.Range("I134").Value = 2910.5652
$ws.Range("J134").Value = 1719
$ws.Range("K134").Value = 8731.695599999999
$ws.Range("L134").Value = 5157
$ws.Range("M134").Value = -6196.695599999999
$ws.Range("N134").Value = -10227

$ws.Range("H136").Value = 1371.1578
$ws.Range("I136").Value = 1012.0833
$ws.Range("K136").Value = 3036.2499
$ws.Range("M136").Value = -486.2498999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 4360937.5
$ws.Range("I33").Value = 7142937
$ws.Range("J33").Value = 33383.332
$ws.Range("K33").Value = 42857622
$ws.Range("L33").Value = 200299.992
$ws.Range("M33").Value = -42857339
$ws.Range("N33").Value = -200865.992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 28781092
$ws.Range("I122").Value = 36719720
$ws.Range("K122").Value = 110159160
$ws.Range("M122").Value = -110156710

$ws.Range("H126").Value = 9115.214
$ws.Range("I126").Value = 12879.111
$ws.Range("J126").Value = 2340.2
$ws.Range("K126").Value = 38637.333
$ws.Range("L126").Value = 7020.599999999999
$ws.Range("M126").Value = -36167.333
$ws.Range("N126").Value = -11960.6

$ws.Range("H132").Value = 2609.2
$ws.Range("I132").Value = 2561.889
$ws.Range("J132").Value = 2698.842
$ws.Range("K132").Value = 7685.667
$ws.Range("L132").Value = 8096.526
$ws.Range("M132").Value = -5155.667
$ws.Range("N132").Value = -13156.526

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 100040400
$ws.Range("I93").Value = 67167.664
$ws.Range("K93").Value = 67167.664
$ws.Range("M93").Value = -65919.664

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H100").Value = 1164.3334
$ws.Range("I100").Value = 996.5
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 996.5
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -455.5
$ws.Range("N100").Value = -2582

$ws.Range("H132").Value = 7395949.5
$ws.Range("I132").Value = 10662950
$ws.Range("K132").Value = 31988850
$ws.Range("M132").Value = -31986320

$ws.Range("H136").Value = 13744.1875
$ws.Range("I136").Value = 19535.285
$ws.Range("J136").Value = 9240
$ws.Range("K136").Value = 58605.855
$ws.Range("L136").Value = 27720
$ws.Range("M136").Value = -56055.855
$ws.Range("N136").Value = -32820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14800.5
$ws.Range("I62").Value = 4750
$ws.Range("K62").Value = 4750
$ws.Range("M62").Value = -4126

$ws.Range("H65").Value = 14800.5
$ws.Range("I65").Value = 4750
$ws.Range("K65").Value = 23750
$ws.Range("M65").Value = -20630

$ws.Range("H107").Value = 40000452
$ws.Range("I107").Value = 50000396
$ws.Range("J107").Value = 678.4
$ws.Range("K107").Value = 150001188
$ws.Range("L107").Value = 2035.2
$ws.Range("M107").Value = -149999268
$ws.Range("N107").Value = -5875.2

$ws.Range("H113").Value = 813.37036
$ws.Range("I113").Value = 765.9375
$ws.Range("J113").Value = 882.36365
$ws.Range("K113").Value = 2297.8125
$ws.Range("L113").Value = 2647.09095
$ws.Range("M113").Value = -127.8125
$ws.Range("N113").Value = -6987.09095

$ws.Range("H132").Value = 1682.5946
$ws.Range("I132").Value = 1045.4286
$ws.Range("J132").Value = 2518.875
$ws.Range("K132").Value = 3136.2858
$ws.Range("L132").Value = 7556.625
$ws.Range("M132").Value = -606.2857999999997
$ws.Range("N132").Value = -12616.625

$ws.Range("H136").Value = 3054.6553
$ws.Range("I136").Value = 3819.375
$ws.Range("J136").Value = 2113.4614
$ws.Range("K136").Value = 11458.125
$ws.Range("L136").Value = 6340.3842
$ws.Range("M136").Value = -8908.125
$ws.Range("N136").Value = -11440.3842
